$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "98.316.10"
$ws.Range("E2").Value = "  +0.72%  "

$ws.Range("D3").Value = "3.511.63"
$ws.Range("E3").Value = "  +4.16%  "

$ws.Range("E4").Value = "  -0.01%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "253.97"
$ws.Range("E5").Value = "  +0.93%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "668.19"
$ws.Range("E6").Value = "  +0.41%  "

$ws.Range("E7").Value = "  +4.40%  "

$ws.Range("E8").Value = "  +1.25%  "

$ws.Range("E9").Value = "  +2.42%  "

$ws.Range("E10").Value = "  -0.01%  "

$ws.Range("D11").Value = "3.509.62"
$ws.Range("E11").Value = "  +4.04%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "46.07"
$ws.Range("E12").Value = "  +11.83%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.211"
$ws.Range("E13").Value = "  -0.39%  "

$ws.Range("D14").Value = "98.349.46"
$ws.Range("E14").Value = "  +1.00%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.22"
$ws.Range("E15").Value = "  +0.97%  "

$ws.Range("D16").Value = "4.186.01"
$ws.Range("E16").Value = "  +4.63%  "

$ws.Range("E17").Value = "  +0.58%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "9.03"
$ws.Range("E18").Value = "  +4.39%  "

$ws.Range("D19").Value = "3.481.46"
$ws.Range("E19").Value = "  +3.45%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "18.88"
$ws.Range("E20").Value = "  +10.69%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.83"
$ws.Range("E21").Value = "  +8.37%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.528"
$ws.Range("E22").Value = "  -8.21%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "524.92"
$ws.Range("E23").Value = "  +4.10%  "

$ws.Range("E24").Value = "  +1.76%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.0000205"
$ws.Range("E25").Value = "  +1.64%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "6.82"
$ws.Range("E26").Value = "  +8.41%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "98.09"
$ws.Range("E27").Value = "  +2.70%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "12.71"
$ws.Range("E28").Value = "  +2.73%  "

$ws.Range("B29").Value = "WrappedeETH"
$ws.Range("C29").Value = "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
$ws.Range("D29").Value = "3.695.50"
$ws.Range("E29").Value = "  +3.35%  "

$ws.Range("B30").Value = "InternetComputer(DFINITY)"
$ws.Range("C30").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "12.48"
$ws.Range("E30").Value = "  +10.57%  "

$ws.Range("B31").Value = "PancakeSwap"
$ws.Range("C31").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.92"
$ws.Range("E31").Value = "  +12.64%  "

$ws.Range("B32").Value = "Hedera"
$ws.Range("C32").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.147"
$ws.Range("E32").Value = "  -2.42%  "

$ws.Range("B33").Value = "Dai"
$ws.Range("C33").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.999"
$ws.Range("E33").Value = "  +0.29%  "

$ws.Range("B34").Value = "Cronos"
$ws.Range("C34").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.191"
$ws.Range("E34").Value = "  +0.14%  "

$ws.Range("B35").Value = "PolygonEcosystemToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.590"
$ws.Range("E35").Value = "  +5.62%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "30.91"
$ws.Range("E36").Value = "  +6.91%  "

$ws.Range("B37").Value = "Binance-PegBSC-USD"
$ws.Range("C37").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.997"
$ws.Range("E37").Value = "  -0.28%  "

$ws.Range("B38").Value = "Fetch.AI"
$ws.Range("C38").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.54"
$ws.Range("E38").Value = "  +3.03%  "

$ws.Range("B39").Value = "RenderToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "8.04"
$ws.Range("E39").Value = "  +2.37%  "

$ws.Range("B40").Value = "Kaspa"
$ws.Range("C40").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.156"
$ws.Range("E40").Value = "  +3.50%  "

$ws.Range("B41").Value = "Bittensor"
$ws.Range("C41").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "532.03"
$ws.Range("E41").Value = "  -0.40%  "

$ws.Range("B42").Value = "USDe"
$ws.Range("C42").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.00"
$ws.Range("E42").Value = "  -0.03%  "

$ws.Range("B43").Value = "ARBITRUM"
$ws.Range("C43").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.917"
$ws.Range("E43").Value = "  +7.74%  "

$ws.Range("B44").Value = "ImmutableX"
$ws.Range("C44").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.79"
$ws.Range("E44").Value = "  +5.87%  "

$ws.Range("B45").Value = "WhiteBITCoin"
$ws.Range("C45").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "24.46"
$ws.Range("E45").Value = "  -0.91%  "

$ws.Range("B46").Value = "VeChain"
$ws.Range("C46").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0434"
$ws.Range("E46").Value = "  +1.30%  "

$ws.Range("B47").Value = "Filecoin"
$ws.Range("C47").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "5.81"
$ws.Range("E47").Value = "  +2.52%  "

$ws.Range("B48").Value = "MantraDAO"
$ws.Range("C48").Value = "https://coinranking.com/coin/cTdD8lD-6+mantradao-om"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.66"
$ws.Range("E48").Value = "  -0.61%  "

$ws.Range("B49").Value = "Cosmos"
$ws.Range("C49").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.68"
$ws.Range("E49").Value = "  -3.26%  "

$ws.Range("B50").Value = "Stacks"
$ws.Range("C50").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.22"
$ws.Range("E50").Value = "  +8.72%  "

$ws.Range("B51").Value = "OKB"
$ws.Range("C51").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "55.32"
$ws.Range("E51").Value = "  +2.06%  "
